$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: add "na" notes in the Cal Cond columns (G,H) ---
$ws.Range("G24").Value = "na"
$ws.Range("H24").Value = "na"

# --- New row 25 ---
$ws.Range("A25").NumberFormat = "d-mmm-yy"
$ws.Range("A25").Value = 44212
$ws.Range("B25").Value = 316
$ws.Range("C25").NumberFormat = "h:mm:ss"
$ws.Range("C25").Value = 0.34375
$ws.Range("D25").NumberFormat = "h:mm:ss"
$ws.Range("D25").Value = 0.34699074074074071
$ws.Range("E25").NumberFormat = "h:mm:ss"
$ws.Range("E25").Font.Bold = $true
$ws.Range("E25").Value = 0.40069444444444446
$ws.Range("F25").NumberFormat = "h:mm:ss"
$ws.Range("F25").Font.Bold = $true
$ws.Range("F25").Value = 0.49907407407407406
$ws.Range("G25").NumberFormat = "h:mm:ss"
$ws.Range("G25").Value = 0.54664351851851845
$ws.Range("H25").NumberFormat = "h:mm:ss"
$ws.Range("H25").Value = 0.55532407407407403
$ws.Range("I25").NumberFormat = "#,##0"
$ws.Range("I25").Value = 50000

# --- New row 26 ---
$ws.Range("A26").NumberFormat = "d-mmm-yy"
$ws.Range("A26").Value = 44212
$ws.Range("B26").Value = 316
$ws.Range("C26").NumberFormat = "h:mm:ss"
$ws.Range("C26").Value = 0.34791666666666665
$ws.Range("D26").NumberFormat = "h:mm:ss"
$ws.Range("D26").Value = 0.35000000000000003
$ws.Range("E26").NumberFormat = "h:mm:ss"
$ws.Range("E26").Font.Bold = $true
$ws.Range("E26").Value = 0.40069444444444446
$ws.Range("F26").NumberFormat = "h:mm:ss"
$ws.Range("F26").Font.Bold = $true
$ws.Range("F26").Value = 0.49907407407407406
$ws.Range("G26").NumberFormat = "h:mm:ss"
$ws.Range("G26").Value = 0.5625
$ws.Range("H26").NumberFormat = "h:mm:ss"
$ws.Range("H26").Value = 0.57777777777777783
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("I26").Value = 1415
$ws.Range("J26").Value = "pre-calibration not stable. Post-calibration took 10 minutes to stabilize"

# --- New row 27 ---
$ws.Range("A27").NumberFormat = "d-mmm-yy"
$ws.Range("A27").Value = 44212
$ws.Range("B27").Value = 354
$ws.Range("C27").NumberFormat = "h:mm:ss"
$ws.Range("C27").Value = 0.3444444444444445
$ws.Range("D27").NumberFormat = "h:mm:ss"
$ws.Range("D27").Value = 0.34791666666666665
$ws.Range("E27").NumberFormat = "h:mm:ss"
$ws.Range("E27").Font.Bold = $true
$ws.Range("E27").Value = 0.40069444444444446
$ws.Range("F27").NumberFormat = "h:mm:ss"
$ws.Range("F27").Font.Bold = $true
$ws.Range("F27").Value = 0.49907407407407406
$ws.Range("G27").NumberFormat = "h:mm:ss"
$ws.Range("G27").Value = 0.54664351851851845
$ws.Range("H27").NumberFormat = "h:mm:ss"
$ws.Range("H27").Value = 0.55532407407407403
$ws.Range("I27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 50000

# --- New row 28 ---
$ws.Range("A28").NumberFormat = "d-mmm-yy"
$ws.Range("A28").Value = 44212
$ws.Range("B28").Value = 354
$ws.Range("C28").NumberFormat = "h:mm:ss"
$ws.Range("C28").Value = 0.35046296296296298
$ws.Range("D28").NumberFormat = "h:mm:ss"
$ws.Range("D28").Value = 0.35185185185185186
$ws.Range("E28").NumberFormat = "h:mm:ss"
$ws.Range("E28").Font.Bold = $true
$ws.Range("E28").Value = 0.40069444444444446
$ws.Range("F28").NumberFormat = "h:mm:ss"
$ws.Range("F28").Font.Bold = $true
$ws.Range("F28").Value = 0.49907407407407406
$ws.Range("G28").NumberFormat = "h:mm:ss"
$ws.Range("G28").Value = 0.5625
$ws.Range("H28").NumberFormat = "h:mm:ss"
$ws.Range("H28").Value = 0.57777777777777783
$ws.Range("I28").NumberFormat = "#,##0"
$ws.Range("I28").Value = 1415

# --- Update the sheet view: scroll position and selection ---
$ws.Activate()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 4
$win.ScrollColumn = 2
$ws.Range("G25:I25").Select()
